# caArray_Risks.xlsx -- "Updated agenda, project plan, action items and risks."
#
# Changes applied:
#   1. Row 11 (R10): Owner corrected from "Rashmi and Shine" to "Eve Shalley".
#   2. Row 13 (R12): Status flipped from Open -> Closed, row re-shaded to match
#      the "closed" look used elsewhere on the sheet (light-gray fill).
#   3. Row 20 (R19): Status flipped from Open -> Closed, re-shaded the same way,
#      and its row height tightened from 45 to 30.
#   4. A brand-new risk, R20, is appended as row 21 (Open status, "no fill"
#      look used for open risks), with Owner/Risk/Resolution text and a 90pt
#      row height to fit the longer mitigation text.
#   5. The view is scrolled/selected to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. R10 owner correction (C11) ---------------------------------------
$ws.Range("C11").Value = "Eve Shalley"

# --- 2. R12 (row 13): mark Closed + re-shade ------------------------------
$ws.Range("E13").Value = "Closed"
$ws.Range("A13:F13").Interior.ColorIndex = 15

# --- 3. R19 (row 20): mark Closed + re-shade + shrink row height ---------
$ws.Range("E20").Value = "Closed"
$ws.Range("A20:F20").Interior.ColorIndex = 15
$ws.Rows(20).RowHeight = 30

# --- 4. New risk R20 (row 21) ---------------------------------------------
# Clone the look of an existing "Open" risk row (row 19) so the new row gets
# the same font/alignment/no-fill formatting without creating extra style
# entries, then fill in the new risk's data.
$ws.Range("A19:F19").Copy() | Out-Null
$ws.Range("A21:F21").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A21").Value = "R20"
$ws.Range("F21").Value = "Updating AHP builds is on hold due to ongoing discussions among Juli, JJ and Doug Hosier about Systems team support. This is a risk since future code changes will happen on GitHub and we will not be able to test on the Dev/QA tiers. (The QA team is scheduled to start regression testing the next release on the QA tier on March 12.)"
$ws.Range("C21").Value = "Juli, JJ and Doug Hosier"
$ws.Range("B21").Value = "After OSDI migration to GitHub, code cannot be tested on NCI tiers due to the AntHill Pro updates being on hold"
$ws.Range("D21").Value = "High"
$ws.Range("E21").Value = "Open"
$ws.Rows(21).RowHeight = 90

# --- 5. Scroll/selection bookkeeping --------------------------------------
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("B22").Select() | Out-Null
